$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Update remark/comment text in column G to change format from "Name：Code" to "Code:Name"
$ws.Range("G10").Value = "1:逾清償1期`n2:逾清償2期`n3:逾清償3-6期"
$ws.Range("G11").Value = "310:短期放款`n320:中期放款`n330:長期放款"
